# Add data for 2025-04-18
# Updates the 2025 (and a handful of prior-year) violent-crime counts
# across the citywide summary, the by-neighborhood rollup, and the
# individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 1726
$ws.Range("K3").Value = 8180
$ws.Range("L3").Value = 1750
$ws.Range("F4").Value = 1928
$ws.Range("G4").Value = 1501
$ws.Range("K4").Value = 1758
$ws.Range("L4").Value = 493
$ws.Range("K6").Value = 9123
$ws.Range("L6").Value = 1638
$ws.Range("F7").Value = 24121
$ws.Range("G7").Value = 24728
$ws.Range("K7").Value = 27548
$ws.Range("L7").Value = 5713

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L2").Value = 4
$ws.Range("L7").Value = 11

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 112
$ws.Range("L7").Value = 344

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 286
$ws.Range("L2").Value = 57
$ws.Range("L6").Value = 88
$ws.Range("K7").Value = 1150
$ws.Range("L7").Value = 250

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 258
$ws.Range("L2").Value = 63
$ws.Range("L3").Value = 58
$ws.Range("L6").Value = 64
$ws.Range("K7").Value = 909
$ws.Range("L7").Value = 205

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 22
$ws.Range("L3").Value = 39
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L5").Value = 19
$ws.Range("L7").Value = 187
$ws.Range("L8").Value = 344
$ws.Range("L11").Value = 103
$ws.Range("K19").Value = 792
$ws.Range("L19").Value = 165
$ws.Range("L20").Value = 155
$ws.Range("L27").Value = 62
$ws.Range("L29").Value = 287
$ws.Range("L30").Value = 28
$ws.Range("K33").Value = 1150
$ws.Range("L33").Value = 250
$ws.Range("L34").Value = 37
$ws.Range("L36").Value = 86
$ws.Range("K37").Value = 909
$ws.Range("L37").Value = 205
$ws.Range("L41").Value = 28
$ws.Range("L42").Value = 179
$ws.Range("L45").Value = 9
$ws.Range("L51").Value = 69
$ws.Range("L52").Value = 121
$ws.Range("F63").Value = 213
$ws.Range("G63").Value = 303
$ws.Range("K63").Value = 88
$ws.Range("L63").Value = 23
$ws.Range("L68").Value = 14
$ws.Range("L69").Value = 11
$ws.Range("L76").Value = 64
$ws.Range("L79").Value = 155
$ws.Range("L85").Value = 296
$ws.Range("L88").Value = 75
$ws.Range("L89").Value = 70
$ws.Range("L90").Value = 56
$ws.Range("L93").Value = 30
$ws.Range("L94").Value = 70
$ws.Range("L97").Value = 56
$ws.Range("L99").Value = 88
$ws.Range("F101").Value = 24121
$ws.Range("G101").Value = 24728
$ws.Range("K101").Value = 27548
$ws.Range("L101").Value = 5713

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 99
$ws.Range("L7").Value = 287

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 234
$ws.Range("L2").Value = 50
$ws.Range("K3").Value = 234
$ws.Range("L4").Value = 8
$ws.Range("K6").Value = 267
$ws.Range("K7").Value = 792
$ws.Range("L7").Value = 165

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 10
$ws.Range("L7").Value = 64

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L4").Value = 1
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 45
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 179

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 56
$ws.Range("L7").Value = 155

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 50
$ws.Range("L6").Value = 47
$ws.Range("L7").Value = 155

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 86

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 58
$ws.Range("L6").Value = 53
$ws.Range("L7").Value = 187

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 30
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 56

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L2").Value = 18
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L3").Value = 19
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 62

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 56

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 20
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 14

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 87
$ws.Range("L3").Value = 125
$ws.Range("L7").Value = 296

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 9

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 40
$ws.Range("L3").Value = 34
$ws.Range("L7").Value = 121
